# Auto-generated update of Leve profit market-price columns (H, I, J, K, L, M, N)
# across all job sheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1919.25
$ws.Range("J17").Value = 2289.8
$ws.Range("L17").Value = 6869.400000000001
$ws.Range("N17").Value = -7205.400000000001
$ws.Range("H18").Value = 1234
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716
$ws.Range("H32").Value = 7039.391
$ws.Range("I32").Value = 10697.2
$ws.Range("J32").Value = 4225.6924
$ws.Range("K32").Value = 10697.2
$ws.Range("L32").Value = 4225.6924
$ws.Range("M32").Value = -10371.2
$ws.Range("N32").Value = -4877.6924
$ws.Range("H51").Value = 5990.1
$ws.Range("J51").Value = 6800.2
$ws.Range("L51").Value = 6800.2
$ws.Range("N51").Value = -7768.2
$ws.Range("H80").Value = 8564.691999999999
$ws.Range("J80").Value = 6661
$ws.Range("L80").Value = 19983
$ws.Range("N80").Value = -21979
$ws.Range("H83").Value = 8564.691999999999
$ws.Range("J83").Value = 6661
$ws.Range("L83").Value = 59949
$ws.Range("N83").Value = -69933
$ws.Range("H97").Value = 796
$ws.Range("J97").Value = 788.25
$ws.Range("L97").Value = 2364.75
$ws.Range("N97").Value = -3356.75
$ws.Range("H98").Value = 1479.75
$ws.Range("I98").Value = 1354.5714
$ws.Range("J98").Value = 2356
$ws.Range("K98").Value = 1354.5714
$ws.Range("L98").Value = 2356
$ws.Range("M98").Value = 143.4286
$ws.Range("N98").Value = -5352
$ws.Range("H122").Value = 1479.75
$ws.Range("I122").Value = 1354.5714
$ws.Range("J122").Value = 2356
$ws.Range("K122").Value = 4063.7142
$ws.Range("L122").Value = 7068
$ws.Range("M122").Value = -1613.7142
$ws.Range("N122").Value = -11968
$ws.Range("H131").Value = 5379.091
$ws.Range("I131").Value = 4917.4
$ws.Range("K131").Value = 14752.2
$ws.Range("M131").Value = -9712.199999999999
$ws.Range("H135").Value = 17543.27
$ws.Range("I135").Value = 1402.9565
$ws.Range("K135").Value = 12626.6085
$ws.Range("M135").Value = -10091.6085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 9692.714
$ws.Range("I5").Value = 9641.5
$ws.Range("K5").Value = 9641.5
$ws.Range("M5").Value = -9529.5
$ws.Range("H63").Value = 4609.8
$ws.Range("I63").Value = 1799
$ws.Range("J63").Value = 5312.5
$ws.Range("K63").Value = 1799
$ws.Range("L63").Value = 5312.5
$ws.Range("M63").Value = -1113
$ws.Range("N63").Value = -6684.5
$ws.Range("H66").Value = 4609.8
$ws.Range("I66").Value = 1799
$ws.Range("J66").Value = 5312.5
$ws.Range("K66").Value = 8995
$ws.Range("L66").Value = 26562.5
$ws.Range("M66").Value = -5563
$ws.Range("N66").Value = -33426.5
$ws.Range("H97").Value = 1507.7333
$ws.Range("I97").Value = 1001.625
$ws.Range("J97").Value = 3532.1667
$ws.Range("K97").Value = 1001.625
$ws.Range("L97").Value = 3532.1667
$ws.Range("M97").Value = -505.625
$ws.Range("N97").Value = -4524.1667
$ws.Range("H132").Value = 2702.7083
$ws.Range("I132").Value = 897.5599999999999
$ws.Range("J132").Value = 4664.826
$ws.Range("K132").Value = 2692.68
$ws.Range("L132").Value = 13994.478
$ws.Range("M132").Value = -162.6799999999998
$ws.Range("N132").Value = -19054.478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 9692.714
$ws.Range("I4").Value = 9641.5
$ws.Range("K4").Value = 9641.5
$ws.Range("M4").Value = -9526.5
$ws.Range("H99").Value = 1994.0667
$ws.Range("I99").Value = 1083.1818
$ws.Range("K99").Value = 1083.1818
$ws.Range("M99").Value = 414.8181999999999
$ws.Range("H107").Value = 103580
$ws.Range("I107").Value = 128225
$ws.Range("K107").Value = 128225
$ws.Range("M107").Value = -126305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1283.1666
$ws.Range("I122").Value = 1283.1666
$ws.Range("K122").Value = 3849.4998
$ws.Range("M122").Value = -1399.4998
$ws.Range("H132").Value = 36865.94
$ws.Range("I132").Value = 41114.8
$ws.Range("K132").Value = 123344.4
$ws.Range("M132").Value = -120814.4
$ws.Range("H134").Value = 2910.3125
$ws.Range("I134").Value = 2588.6155
$ws.Range("K134").Value = 7765.8465
$ws.Range("M134").Value = -5230.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 868334.3
$ws.Range("I2").Value = 1389305.8
$ws.Range("K2").Value = 8335834.800000001
$ws.Range("M2").Value = -8335721.800000001
$ws.Range("H116").Value = 7938.625
$ws.Range("I116").Value = 2887.8333
$ws.Range("J116").Value = 9622.223
$ws.Range("K116").Value = 8663.499899999999
$ws.Range("L116").Value = 28866.669
$ws.Range("M116").Value = -5221.499899999999
$ws.Range("N116").Value = -35750.669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10346
$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("N30").Value = -10210
$ws.Range("H70").Value = 8017.9
$ws.Range("I70").Value = 8059.875
$ws.Range("K70").Value = 8059.875
$ws.Range("M70").Value = -7789.875
$ws.Range("H73").Value = 8017.9
$ws.Range("I73").Value = 8059.875
$ws.Range("K73").Value = 8059.875
$ws.Range("M73").Value = -7123.875
$ws.Range("H122").Value = 3108.1428
$ws.Range("I122").Value = 2622.9285
$ws.Range("J122").Value = 4078.5715
$ws.Range("K122").Value = 7868.7855
$ws.Range("L122").Value = 12235.7145
$ws.Range("M122").Value = -5418.7855
$ws.Range("N122").Value = -17135.7145
$ws.Range("H132").Value = 5051.7334
$ws.Range("I132").Value = 4752.077
$ws.Range("K132").Value = 14256.231
$ws.Range("M132").Value = -11726.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2070.8
$ws.Range("I22").Value = 2001
$ws.Range("K22").Value = 2001
$ws.Range("M22").Value = -1706
$ws.Range("H25").Value = 55833.332
$ws.Range("I25").Value = 53333.332
$ws.Range("J25").Value = 56666.668
$ws.Range("K25").Value = 53333.332
$ws.Range("L25").Value = 56666.668
$ws.Range("M25").Value = -53103.332
$ws.Range("N25").Value = -57126.668
$ws.Range("H27").Value = 2070.8
$ws.Range("I27").Value = 2001
$ws.Range("K27").Value = 2001
$ws.Range("M27").Value = -1894
$ws.Range("H61").Value = 1003.4
$ws.Range("I61").Value = 1009.8571
$ws.Range("J61").Value = 988.3333
$ws.Range("K61").Value = 1009.8571
$ws.Range("L61").Value = 988.3333
$ws.Range("M61").Value = -807.8570999999999
$ws.Range("N61").Value = -1392.3333
$ws.Range("H113").Value = 1003.4
$ws.Range("I113").Value = 1009.8571
$ws.Range("J113").Value = 988.3333
$ws.Range("K113").Value = 1009.8571
$ws.Range("L113").Value = 988.3333
$ws.Range("M113").Value = 1160.1429
$ws.Range("N113").Value = -5328.3333
$ws.Range("H122").Value = 8942.483
$ws.Range("I122").Value = 9847.888999999999
$ws.Range("K122").Value = 29543.667
$ws.Range("M122").Value = -27093.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 26679.334
$ws.Range("I18").Value = 20000
$ws.Range("K18").Value = 20000
$ws.Range("M18").Value = -19827
$ws.Range("H100").Value = 1600
$ws.Range("I100").Value = 950
$ws.Range("K100").Value = 1900
$ws.Range("M100").Value = -1359
$ws.Range("H113").Value = 1280.7142
$ws.Range("I113").Value = 1322.8
$ws.Range("J113").Value = 1257.3334
$ws.Range("K113").Value = 3968.4
$ws.Range("L113").Value = 3772.0002
$ws.Range("M113").Value = -1798.4
$ws.Range("N113").Value = -8112.0002
$ws.Range("H122").Value = 65576.11
$ws.Range("I122").Value = 80529.73
$ws.Range("K122").Value = 241589.19
$ws.Range("M122").Value = -239139.19
$ws.Range("H136").Value = 22336.062
$ws.Range("I136").Value = 24495.896
$ws.Range("J136").Value = 1457.6666
$ws.Range("K136").Value = 73487.68799999999
$ws.Range("L136").Value = 4372.9998
$ws.Range("M136").Value = -70937.68799999999
$ws.Range("N136").Value = -9472.9998
